$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 13079.5
$ws.Range("J88").Value = 17995.75
$ws.Range("L88").Value = 17995.75
$ws.Range("N88").Value = -18807.75

$ws.Range("H91").Value = 13079.5
$ws.Range("J91").Value = 17995.75
$ws.Range("L91").Value = 17995.75
$ws.Range("N91").Value = -20803.75

$ws.Range("H135").Value = 1051.9333
$ws.Range("J135").Value = 1706
$ws.Range("L135").Value = 15354
$ws.Range("N135").Value = -20424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 2537.5
$ws.Range("I39").Value = 975
$ws.Range("J39").Value = 4100
$ws.Range("K39").Value = 975
$ws.Range("L39").Value = 4100
$ws.Range("M39").Value = -455
$ws.Range("N39").Value = -5140

$ws.Range("H42").Value = 18765.25
$ws.Range("J42").Value = 18765.25
$ws.Range("L42").Value = 18765.25
$ws.Range("N42").Value = -19737.25

$ws.Range("H61").Value = 3626.3225
$ws.Range("I61").Value = 2133.963
$ws.Range("K61").Value = 2133.963
$ws.Range("M61").Value = -1921.963

$ws.Range("H122").Value = 4838.6665
$ws.Range("I122").Value = 4004.5
$ws.Range("K122").Value = 12013.5
$ws.Range("M122").Value = -9563.5

$ws.Range("H132").Value = 2869.9614
$ws.Range("I132").Value = 2294.55
$ws.Range("K132").Value = 6883.650000000001
$ws.Range("M132").Value = -4353.650000000001

$ws.Range("H136").Value = 3626.3225
$ws.Range("I136").Value = 2133.963
$ws.Range("K136").Value = 6401.889000000001
$ws.Range("M136").Value = -3851.889000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8179.85
$ws.Range("I105").Value = 4007.9092
$ws.Range("K105").Value = 4007.9092
$ws.Range("M105").Value = -2260.9092

$ws.Range("H107").Value = 1360.7
$ws.Range("I107").Value = 1400.8889
$ws.Range("K107").Value = 1400.8889
$ws.Range("M107").Value = 519.1111000000001

$ws.Range("H134").Value = 1930.125
$ws.Range("J134").Value = 12506
$ws.Range("L134").Value = 37518
$ws.Range("N134").Value = -42588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28627.195
$ws.Range("I31").Value = 2911.423
$ws.Range("J31").Value = 73201.2
$ws.Range("K31").Value = 2911.423
$ws.Range("L31").Value = 73201.2
$ws.Range("M31").Value = -2616.423
$ws.Range("N31").Value = -73791.2

$ws.Range("H34").Value = 28627.195
$ws.Range("I34").Value = 2911.423
$ws.Range("J34").Value = 73201.2
$ws.Range("K34").Value = 2911.423
$ws.Range("L34").Value = 73201.2
$ws.Range("M34").Value = -2709.423
$ws.Range("N34").Value = -73605.2

$ws.Range("H58").Value = 3657.1538
$ws.Range("I58").Value = 1733.625
$ws.Range("K58").Value = 1733.625
$ws.Range("M58").Value = -1530.625

$ws.Range("H107").Value = 803.7742
$ws.Range("I107").Value = 703.95
$ws.Range("J107").Value = 985.2727
$ws.Range("K107").Value = 703.95
$ws.Range("L107").Value = 985.2727
$ws.Range("M107").Value = 1216.05
$ws.Range("N107").Value = -4825.2727

$ws.Range("H134").Value = 2794.35
$ws.Range("I134").Value = 1791.9412
$ws.Range("K134").Value = 5375.8236
$ws.Range("M134").Value = -2840.8236

$ws.Range("H136").Value = 3657.1538
$ws.Range("I136").Value = 1733.625
$ws.Range("K136").Value = 5200.875
$ws.Range("M136").Value = -2650.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1187.0769
$ws.Range("I113").Value = 760.625
$ws.Range("J113").Value = 1869.4
$ws.Range("K113").Value = 2281.875
$ws.Range("L113").Value = 5608.200000000001
$ws.Range("M113").Value = -111.875
$ws.Range("N113").Value = -9948.200000000001

$ws.Range("H121").Value = 2949.5833
$ws.Range("I121").Value = 4000
$ws.Range("J121").Value = 2739.5
$ws.Range("K121").Value = 12000
$ws.Range("L121").Value = 8218.5
$ws.Range("M121").Value = -10690
$ws.Range("N121").Value = -10838.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 833666.7
$ws.Range("J12").Value = 5000000
$ws.Range("L12").Value = 5000000
$ws.Range("N12").Value = -5000280

$ws.Range("H17").Value = 5303.3335
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 6164
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 6164
$ws.Range("M17").Value = -832
$ws.Range("N17").Value = -6500

$ws.Range("M23").ClearContents()
$ws.Range("H23").Value = 2500
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 2500
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 2500
$ws.Range("N23").Value = -2946

$ws.Range("H24").Value = 22183.8
$ws.Range("I24").Value = 23970.666
$ws.Range("K24").Value = 23970.666
$ws.Range("M24").Value = -23797.666

$ws.Range("H80").Value = 6329.8823
$ws.Range("I80").Value = 4239.8
$ws.Range("J80").Value = 7200.75
$ws.Range("K80").Value = 4239.8
$ws.Range("L80").Value = 7200.75
$ws.Range("M80").Value = -3241.8
$ws.Range("N80").Value = -9196.75

$ws.Range("H83").Value = 6329.8823
$ws.Range("I83").Value = 4239.8
$ws.Range("J83").Value = 7200.75
$ws.Range("K83").Value = 21199
$ws.Range("L83").Value = 36003.75
$ws.Range("M83").Value = -16207
$ws.Range("N83").Value = -45987.75

$ws.Range("H111").Value = 33315.5
$ws.Range("J111").Value = 52631
$ws.Range("L111").Value = 52631
$ws.Range("N111").Value = -58765

$ws.Range("H132").Value = 33007
$ws.Range("I132").Value = 45743.22
$ws.Range("J132").Value = 6376.727
$ws.Range("K132").Value = 137229.66
$ws.Range("L132").Value = 19130.181
$ws.Range("M132").Value = -134699.66
$ws.Range("N132").Value = -24190.181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3414.1667
$ws.Range("I46").Value = 1659
$ws.Range("K46").Value = 1659
$ws.Range("M46").Value = -1471

$ws.Range("N54").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0

$ws.Range("H132").Value = 2942.3953
$ws.Range("I132").Value = 1778.2333
$ws.Range("K132").Value = 5334.699900000001
$ws.Range("M132").Value = -2804.699900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6980.6
$ws.Range("J62").Value = 4999.5
$ws.Range("L62").Value = 4999.5
$ws.Range("N62").Value = -6247.5

$ws.Range("H65").Value = 6980.6
$ws.Range("J65").Value = 4999.5
$ws.Range("L65").Value = 24997.5
$ws.Range("N65").Value = -31237.5

$ws.Range("H132").Value = 4182.7646
$ws.Range("I132").Value = 3558.5417
$ws.Range("K132").Value = 10675.6251
$ws.Range("M132").Value = -8145.625100000001
